$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "La Laja"
$ws.Range("C2").Value = "Carboneras"
$ws.Range("B3").Value = "Cascabeles"
$ws.Range("C3").Value = "Jazmin"
$ws.Range("B4").Value = "Giro"
$ws.Range("C4").Value = "Alonsos"

$ws.Range("D15").Select()
